# Add 2022-Q4 data:
#  - "总计" gains a new summary row for 2022-Q4 (old 2022-Q3 row shifts down).
#  - A new "2022-Q4" sheet is inserted (holding per-fund data), the existing
#    "2022-Q3" per-fund sheet is preserved unchanged but moves to 3rd place.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q3" fund-detail sheet so that an unmodified
#    copy of it ends up in 3rd position (same data, same formatting).
# ---------------------------------------------------------------------------
$q3.Copy($null, $q3)
$q3copy = $wb.Worksheets.Item(3)

# Free up the "2022-Q3" name on the original sheet before claiming it on the
# copy (Excel won't let two sheets share a name).
$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Turn the original sheet (now named "2022-Q4", 2nd position) into the new
#    quarter's fund-detail table.
# ---------------------------------------------------------------------------
$q4 = $q3
$q4.Cells.Clear()

# Copy the header-row style (bold + border + centered) already used on the
# "总计" sheet so formatting matches what the workbook already uses.
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns D:G hold decimal-looking numbers that must stay TEXT (matching the
# existing "2022-Q3" sheet), and column B holds fund codes that must keep
# any leading zeros. A leading apostrophe forces text entry; re-applying the
# built-in "Normal" style afterwards drops the transient quote-prefix flag
# that Excel adds so the cell format stays plain (same as its neighbours).
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'161611"
$q4.Range("C2").Value = "融通内需驱动混合A/B"
$q4.Range("D2").Value = "'9.28"
$q4.Range("E2").Value = "'92.61"
$q4.Range("F2").Value = "'5.36"
$q4.Range("G2").Value = "'0.4974"
$q4.Range("H2").Value = 3

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'014109"
$q4.Range("C3").Value = "融通内需驱动混合C"
$q4.Range("D3").Value = "'3.63"
$q4.Range("E3").Value = "'92.61"
$q4.Range("F3").Value = "'5.36"
$q4.Range("G3").Value = "'0.1946"
$q4.Range("H3").Value = 3

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'014106"
$q4.Range("C4").Value = "融通成长30灵活配置混合C"
$q4.Range("D4").Value = "'2.19"
$q4.Range("E4").Value = "'94.02"
$q4.Range("F4").Value = "'4.87"
$q4.Range("G4").Value = "'0.1067"
$q4.Range("H4").Value = 4

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'002252"
$q4.Range("C5").Value = "融通成长30灵活配置混合A/B"
$q4.Range("D5").Value = "'1.65"
$q4.Range("E5").Value = "'94.02"
$q4.Range("F5").Value = "'4.87"
$q4.Range("G5").Value = "'0.0804"
$q4.Range("H5").Value = 4

$q4.Range("B2:B5").Style = "Normal"
$q4.Range("D2:G5").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert the new 2022-Q4 row and push the
#    2022-Q3 row down to row 3 (values unchanged).
# ---------------------------------------------------------------------------
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.47

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.88
